$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Select() | Out-Null

# Insert a new blank column before column N (14th column), pushing the
# existing "Late"/"heading"/"Outstanding" columns one place to the right.
$col = $ws.Columns.Item(14)
$col.Insert()

# Match the width of the newly inserted column to column M (the column
# immediately to its left), same as the new column picked up in the source.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Leave the selection where the author left it when they saved the file.
$ws.Range("L18").Select() | Out-Null
